$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string in A1
$ws.Range("A1").Value = 'Datos actualizados a 22 de Julio de 2020 a las 15:33'

# Refreshed country data (sorted descending by Casos totales), rows 4-219
$data = @(
  @('Estados Unidos', 4030200, 1631, 1887486, 1997731, 0, 30, 144983),
  @('Brasil', 2166532, 0, 1465970, 618965, 0, 0, 81597),
  @('India', 1201727, 7642, 759168, 413713, 0, 76, 28846),
  @('Rusia', 789190, 5862, 572053, 204392, 0, 165, 12745),
  @('Sudafrica', 381798, 0, 208144, 168286, 0, 0, 5368),
  @('Peru', 362087, 0, 248746, 99762, 0, 0, 13579),
  @('Mexico', 356255, 6859, 227165, 88690, 0, 915, 40400),
  @('Chile', 334683, 0, 306816, 19190, 0, 0, 8677),
  @('España', 313274, 0, 0, 0, 0, 0, 28424),
  @('Reino Unido', 295817, 0, 0, 0, 0, 0, 45422),
  @('Iran', 281413, 2586, 244840, 21720, 0, 219, 14853),
  @('Pakistan', 267428, 1332, 210468, 51283, 0, 38, 5677),
  @('Arabia Saudita', 258156, 2331, 210398, 45157, 0, 44, 2601),
  @('Italia', 244752, 0, 197431, 12248, 0, 0, 35073),
  @('Turquia', 221500, 0, 204011, 11963, 0, 0, 5526),
  @('Banglades', 213254, 2744, 117202, 93301, 0, 42, 2751),
  @('Colombia', 211038, 0, 98840, 105032, 0, 0, 7166),
  @('Alemania', 203890, 0, 188100, 6610, 0, 0, 9180),
  @('Francia', 177338, 0, 79734, 67439, 0, 0, 30165),
  @('Argentina', 136118, 0, 60531, 73081, 0, 16, 2506),
  @('Canada', 111697, 0, 97757, 5078, 0, 0, 8862),
  @('Catar', 107871, 441, 104641, 3067, 0, 3, 163),
  @('Irak', 97159, 0, 64950, 28259, 0, 0, 3950),
  @('Indonesia', 91751, 1882, 50255, 37037, 0, 139, 4459),
  @('Egipto', 89078, 0, 29473, 55206, 0, 0, 4399),
  @('China', 83707, 14, 78840, 233, 0, 0, 4634),
  @('Suecia', 78166, 0, 0, 0, 0, 0, 5646),
  @('Ecuador', 76217, 0, 32725, 38126, 0, 0, 5366),
  @('Kazajistan', 75153, 1685, 46790, 27778, 0, 0, 585),
  @('Filipinas', 72269, 1594, 23623, 46803, 0, 6, 1843),
  @('Oman', 71547, 1660, 47922, 23276, 0, 12, 349),
  @('Bielorrusia', 66521, 173, 59061, 6947, 0, 6, 513),
  @('Belgica', 64258, 164, 17330, 37123, 0, 0, 9805),
  @('Bolivia', 62357, 1366, 19290, 40794, 0, 55, 2273),
  @('Kuwait', 61185, 751, 51520, 9248, 0, 5, 417),
  @('Ucrania', 60995, 829, 33172, 26289, 0, 16, 1534),
  @('Emiratos Arabes Unidos', 57498, 0, 49964, 7193, 0, 0, 341),
  @('Panama', 55153, 0, 30075, 23919, 0, 0, 1159),
  @('Republica Dominicana', 54797, 0, 25976, 27822, 0, 0, 999),
  @('Israel', 54663, 621, 22920, 31313, 0, 5, 430),
  @('Paises Bajos', 52242, 169, 0, 0, 0, 3, 6139),
  @('Portugal', 48898, 0, 33769, 13432, 0, 0, 1697),
  @('Singapur', 48744, 310, 44584, 4133, 0, 0, 27),
  @('Polonia', 41162, 380, 31139, 8381, 0, 6, 1642),
  @('Guatemala', 40229, 0, 26685, 12013, 0, 0, 1531),
  @('Rumania', 40163, 1030, 24663, 13399, 0, 27, 2101),
  @('Nigeria', 37801, 0, 15677, 21319, 0, 0, 805),
  @('Barein', 37316, 0, 33455, 3731, 0, 1, 130),
  @('Afganistan', 35727, 112, 23924, 10613, 0, 4, 1190),
  @('Armenia', 35693, 439, 24766, 10249, 0, 16, 678),
  @('Honduras', 35345, 734, 4144, 30213, 0, 53, 988),
  @('Suiza', 33883, 141, 30500, 1411, 0, 0, 1972),
  @('Ghana', 28989, 0, 25331, 3505, 0, 0, 153),
  @('Kirguistan', 28980, 729, 15536, 12333, 0, 32, 1111),
  @('Azerbaiyan', 28633, 391, 20443, 7805, 0, 9, 385),
  @('Irlanda', 25802, 0, 23364, 685, 0, 0, 1753),
  @('Japon', 25736, 0, 20155, 4593, 0, 0, 988),
  @('Argelia', 24278, 0, 16646, 6532, 0, 0, 1100),
  @('Serbia', 22031, 426, 14047, 7485, 0, 8, 499),
  @('Moldavia', 21442, 0, 14856, 5878, 0, 1, 708),
  @('Austria', 19929, 102, 17849, 1369, 0, 1, 711),
  @('Uzbekistan', 18324, 443, 9554, 8673, 0, 2, 97),
  @('Nepal', 18094, 100, 12684, 5368, 0, 2, 42),
  @('Marruecos', 17742, 0, 15389, 2073, 0, 0, 280),
  @('Camerun', 16522, 0, 13728, 2412, 0, 0, 382),
  @('Kenia', 14805, 637, 6757, 7788, 0, 10, 260),
  @('Costa de Marfil', 14531, 0, 8857, 5581, 0, 0, 93),
  @('Chequia', 14324, 0, 8918, 5046, 0, 0, 360),
  @('Corea del Sur', 13879, 63, 12698, 884, 0, 1, 297),
  @('Dinamarca', 13350, 48, 12274, 465, 0, 0, 611),
  @('El Salvador', 12975, 393, 7119, 5493, 0, 11, 363),
  @('Australia', 12894, 466, 8545, 4221, 0, 2, 128),
  @('Venezuela', 12774, 0, 6983, 5671, 0, 0, 120),
  @('Costa Rica', 11811, 0, 3194, 8549, 0, 0, 68),
  @('Etiopia', 11072, 0, 5448, 5444, 0, 0, 180),
  @('Sudan', 10992, 0, 5707, 4592, 0, 0, 693),
  @('Republica de Macedonia', 9412, 0, 4940, 4040, 0, 0, 432),
  @('Estado de Palestina', 9398, 170, 1932, 7402, 0, 0, 64),
  @('Bulgaria', 9254, 0, 4521, 4420, 0, 0, 313),
  @('Senegal', 9121, 136, 6108, 2836, 0, 3, 177),
  @('Bosnia y Herzegovina', 9115, 328, 4159, 4688, 0, 4, 268),
  @('Noruega', 9059, 6, 8138, 666, 0, 0, 255),
  @('Malasia', 8831, 16, 8566, 142, 0, 0, 123),
  @('Consejo Danes para los Refugiados', 8626, 92, 4790, 3639, 0, 1, 197),
  @('Madagascar', 8162, 614, 4662, 3431, 0, 4, 69),
  @('Finlandia', 7362, 11, 6920, 114, 0, 0, 328),
  @('Haiti', 7146, 46, 4095, 2897, 0, 3, 154),
  @('Tayikistan', 6967, 0, 5683, 1227, 0, 0, 57),
  @('Guayana Francesa', 6851, 0, 4996, 1816, 0, 0, 39),
  @('Guinea', 6652, 0, 5771, 840, 0, 0, 41),
  @('Gabon', 6433, 0, 4034, 2353, 0, 0, 46),
  @('Mauritania', 5985, 0, 3826, 2004, 0, 0, 155),
  @('Luxemburgo', 5725, 0, 4479, 1135, 0, 0, 111),
  @('Republica de Yibuti', 5030, 3, 4904, 68, 0, 2, 58),
  @('Republica de Africa Central', 4561, 0, 1411, 3095, 0, 0, 55),
  @('Croacia', 4530, 108, 3278, 1127, 0, 2, 125),
  @('Hungria', 4366, 19, 3283, 487, 0, 0, 596),
  @('Albania', 4290, 0, 2397, 1776, 0, 0, 117),
  @('Grecia', 4048, 0, 1374, 2477, 0, 0, 197),
  @('Paraguay', 3817, 0, 2307, 1475, 0, 0, 35),
  @('Zambia', 3583, 197, 1677, 1778, 0, 0, 128),
  @('Nicaragua', 3439, 0, 2492, 839, 0, 0, 108),
  @('Tailandia', 3261, 6, 3105, 98, 0, 0, 58),
  @('Malaui', 3149, 104, 1256, 1822, 0, 7, 71),
  @('Somalia', 3135, 0, 1464, 1578, 0, 0, 93),
  @('Guinea Ecuatorial', 3071, 0, 842, 2178, 0, 0, 51),
  @('Maldivas', 3044, 0, 2397, 632, 0, 0, 15),
  @('Libano', 2980, 0, 1577, 1362, 0, 0, 41),
  @('Congo', 2851, 0, 666, 2135, 0, 0, 50),
  @('Mayotte', 2824, 0, 2614, 172, 0, 0, 38),
  @('Sri Lanka', 2730, 0, 2064, 655, 0, 0, 11),
  @('Mali', 2477, 0, 1856, 499, 0, 0, 122),
  @('Cuba', 2449, 0, 2321, 41, 0, 0, 87),
  @('Montenegro', 2381, 0, 426, 1920, 0, 0, 35),
  @('Sudan del Sur', 2211, 0, 1175, 991, 0, 0, 45),
  @('Hong Kong', 2132, 112, 1344, 774, 0, 0, 14),
  @('Cabo Verde', 2107, 0, 1100, 986, 0, 0, 21),
  @('Libia', 2088, 0, 479, 1559, 0, 0, 50),
  @('Eslovaquia', 2058, 37, 1556, 474, 0, 0, 28),
  @('Estonia', 2025, 3, 1912, 44, 0, 0, 69),
  @('Eslovenia', 2006, 29, 1648, 243, 0, 1, 115),
  @('Guinea-Bisau', 1954, 0, 803, 1125, 0, 0, 26),
  @('Lituania', 1951, 2, 1607, 264, 0, 0, 80),
  @('Suazilandia', 1894, 0, 855, 1015, 0, 0, 24),
  @('Islandia', 1840, 1, 1822, 8, 0, 0, 10),
  @('Zimbabue', 1820, 0, 488, 1306, 0, 0, 26),
  @('Sierra Leona', 1727, 0, 1273, 388, 0, 0, 66),
  @('Ruanda', 1655, 0, 848, 802, 0, 0, 5),
  @('Yemen', 1629, 0, 741, 432, 0, 0, 456),
  @('Benin', 1602, 0, 782, 789, 0, 0, 31),
  @('Nueva Zelanda', 1555, 0, 1506, 27, 0, 0, 22),
  @('Mozambique', 1536, 0, 506, 1019, 0, 0, 11),
  @('Tunez', 1389, 0, 1103, 236, 0, 0, 50),
  @('Namibia', 1366, 0, 57, 1302, 0, 0, 7),
  @('Letonia', 1197, 4, 1045, 121, 0, 0, 31),
  @('Surinam', 1131, 0, 705, 405, 0, 0, 21),
  @('Jordania', 1113, 0, 1034, 68, 0, 0, 11),
  @('Niger', 1113, 0, 1018, 26, 0, 0, 69),
  @('Liberia', 1108, 0, 575, 463, 0, 0, 70),
  @('Uruguay', 1096, 0, 929, 134, 0, 0, 33),
  @('Uganda', 1075, 0, 958, 117, 0, 0, 0),
  @('Georgia', 1073, 24, 907, 150, 0, 0, 16),
  @('Burkina Faso', 1066, 1, 917, 96, 0, 0, 53),
  @('Republica de Chipre', 1040, 0, 847, 174, 0, 0, 19),
  @('Republica del Chad', 889, 0, 803, 11, 0, 0, 75),
  @('Principado de Andorra', 884, 0, 803, 29, 0, 0, 52),
  @('Jamaica', 810, 1, 709, 91, 0, 0, 10),
  @('Togo', 790, 0, 560, 215, 0, 0, 15),
  @('Angola', 779, 0, 221, 528, 0, 0, 30),
  @('Santo Tome y Principe', 746, 0, 588, 144, 0, 0, 14),
  @('Crucero', 712, 0, 651, 48, 0, 0, 13),
  @('San Marino', 699, 0, 656, 1, 0, 0, 42),
  @('Malta', 679, 2, 665, 5, 0, 0, 9),
  @('Reunion', 639, 0, 472, 164, 0, 0, 3),
  @('Siria', 540, 0, 160, 349, 0, 0, 31),
  @('Botsuana', 522, 0, 48, 473, 0, 0, 1),
  @('Tanzania', 509, 0, 183, 305, 0, 0, 21),
  @('Taiwan', 455, 0, 440, 8, 0, 0, 7),
  @('Vietnam', 408, 7, 365, 43, 0, 0, 0),
  @('Lesoto', 359, 0, 69, 284, 0, 0, 6),
  @('Mauricio', 343, 0, 332, 1, 0, 0, 10),
  @('Birmania', 341, 0, 280, 55, 0, 0, 6),
  @('Guyana', 339, 0, 163, 157, 0, 0, 19),
  @('Comoras', 337, 0, 319, 11, 0, 0, 7),
  @('Isla de Man', 336, 0, 312, 0, 0, 0, 24),
  @('Burundi', 328, 0, 207, 120, 0, 0, 1),
  @('Mongolia', 287, 0, 213, 74, 0, 0, 0),
  @('Martinica', 262, 0, 98, 149, 0, 0, 15),
  @('Eritrea', 251, 0, 155, 96, 0, 0, 0),
  @('Islas Caimanes', 203, 0, 201, 1, 0, 0, 1),
  @('Camboya', 197, 0, 140, 57, 0, 0, 0),
  @('Guadalupe', 195, 0, 172, 9, 0, 0, 14),
  @('Bahamas', 194, 0, 91, 92, 0, 0, 11),
  @('Islas Feroe', 191, 0, 188, 3, 0, 0, 0),
  @('Gibraltar', 182, 2, 180, 2, 0, 0, 0),
  @('Bermudas', 153, 0, 137, 7, 0, 0, 9),
  @('Brunei', 141, 0, 138, 0, 0, 0, 3),
  @('Trinidad yTobago', 139, 0, 124, 7, 0, 0, 8),
  @('Gambia', 132, 20, 57, 70, 0, 1, 5),
  @('Aruba', 117, 0, 100, 14, 0, 0, 3),
  @('Monaco', 111, 0, 100, 7, 0, 0, 4),
  @('Seychelles', 108, 0, 27, 81, 0, 0, 0),
  @('Barbados', 106, 0, 91, 8, 0, 0, 7),
  @('Butan', 92, 2, 83, 9, 0, 0, 0),
  @('Liechtenstein', 86, 0, 81, 4, 0, 0, 1),
  @('Islas Turcas y Caicos', 81, 0, 22, 57, 0, 0, 2),
  @('San Martin (Parte Holandesa)', 79, 0, 63, 1, 0, 0, 15),
  @('Antigua y Barbuda', 76, 0, 57, 16, 0, 0, 3),
  @('Polinesia Francesa', 62, 0, 60, 2, 0, 0, 0),
  @('San Vicente y las Granadinas', 50, 0, 29, 21, 0, 0, 0),
  @('San Martin (Parte Francesa)', 46, 0, 39, 4, 0, 0, 3),
  @('Macao', 46, 0, 46, 0, 0, 0, 0),
  @('Belice', 43, 1, 23, 18, 0, 0, 2),
  @('Puerto Rico', 39, 0, 1, 36, 0, 0, 2),
  @('Guam', 32, 0, 0, 31, 0, 0, 1),
  @('Curazao', 28, 0, 24, 3, 0, 0, 1),
  @('Papua Nueva Guinea', 27, 0, 8, 19, 0, 0, 0),
  @('Fiyi', 27, 0, 18, 9, 0, 0, 0),
  @('Timor Oriental', 24, 0, 24, 0, 0, 0, 0),
  @('Santa Lucia', 23, 0, 19, 4, 0, 0, 0),
  @('Granada', 23, 0, 23, 0, 0, 0, 0),
  @('Nueva Caledonia', 22, 0, 21, 1, 0, 0, 0),
  @('Laos', 19, 0, 19, 0, 0, 0, 0),
  @('Dominica', 18, 0, 18, 0, 0, 0, 0),
  @('Islas Virgenes de los Estados Unidos', 17, 0, 0, 17, 0, 0, 0),
  @('San Cristobal y Nieves', 17, 0, 15, 2, 0, 0, 0),
  @('Groenlandia', 13, 0, 13, 0, 0, 0, 0),
  @('Islas Malvinas', 13, 0, 13, 0, 0, 0, 0),
  @('Montserrat', 12, 0, 10, 1, 0, 0, 1),
  @('Santa Sede', 12, 0, 12, 0, 0, 0, 0),
  @('Bonaire, San Eustaquio y Saba', 10, 0, 7, 3, 0, 0, 0),
  @('Sahara Occidental', 10, 0, 8, 1, 0, 0, 1),
  @('Islas Virgenes Britanicas', 8, 0, 7, 0, 0, 0, 1),
  @('San Bartolome', 6, 0, 6, 0, 0, 0, 0),
  @('San Pedro y Miquelon', 4, 0, 1, 3, 0, 0, 0),
  @('Anguila', 3, 0, 3, 0, 0, 0, 0)
)

$startRow = 4
for ($i = 0; $i -lt $data.Count; $i++) {
  $row = $data[$i]
  $r = $startRow + $i
  for ($j = 0; $j -lt $row.Count; $j++) {
    $ws.Cells.Item($r, $j + 1).Value = $row[$j]
  }
}
